# Creating PowerShell Modules - step 3 edit
# - Insert two new slides ("Assumptions", "Case") after the title slide
# - Rework the "Problem" slide's bullet content and reflow its diagram shapes

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. New slide 2: "Assumptions"
# ---------------------------------------------------------------------------
$sAssumptions = $p.Slides.Add(2, 2)   # ppLayoutText ("Title and Content")

$sAssumptions.Shapes.Item(1).TextFrame.TextRange.Text = "Assumptions"

$assumptionsBody = $sAssumptions.Shapes.Item(2).TextFrame.TextRange
$assumptionsBody.Text = (
    "We want our code to be easy to maintain`r" +
    "We want our code to be easy to read`r" +
    "We want our code to be easy to refactor or rework`r" +
    "We want to stimulate reusability of code"
)

# ---------------------------------------------------------------------------
# 2. New slide 3: "Case"
# ---------------------------------------------------------------------------
$sCase = $p.Slides.Add(3, 2)          # ppLayoutText ("Title and Content")

$sCase.Shapes.Item(1).TextFrame.TextRange.Text = "Case`t"

$apos = [char]0x2019
$caseBody = $sCase.Shapes.Item(2).TextFrame.TextRange
$caseBody.Text = (
    "We${apos}ve created a solution (repo 1) for a client that includes PowerShell scripts for deploying the solution.`r" +
    "The client asked us to create a second solution (repo 2) for them. `r" +
    "We want to reuse the PowerShell scripts in this second solution."
)

# ---------------------------------------------------------------------------
# 3. Existing "Problem" slide (now at logical position 4): rewrite bullets
#    and reflow the diagram (pictures + red rectangles) further up the slide.
# ---------------------------------------------------------------------------
$sProblem = $p.Slides.Item(4)

$problemBody = $sProblem.Shapes.Item(2).TextFrame.TextRange
$problemBody.Text = (
    "Code sharing across repos`r" +
    "Functions often get copied across, creating a hard-to-debug and hard-to-maintain landscape of repos`r" +
    "`r" +
    "Code sharing within a large repo`r" +
    "Functions provide benefits, but complicate refactoring when using dot sourcing`r" +
    "Nobody likes long relative paths when using dot sourcing`r"
)

# Re-indent the "sub-bullet" lines (2nd, 3rd, 5th and 6th paragraphs) to
# outline level 1 (IndentLevel 2 == OOXML lvl="1"), matching the original
# layout of the slide.
$problemBody.Paragraphs(2,1).IndentLevel = 2
$problemBody.Paragraphs(3,1).IndentLevel = 2
$problemBody.Paragraphs(5,1).IndentLevel = 2
$problemBody.Paragraphs(6,1).IndentLevel = 2

# Move the diagram shapes (3 pictures + 4 red outline rectangles) up the slide.
function Set-ShapePos($shape, $left, $top) {
    $shape.Left = $left / 12700.0
    $shape.Top  = $top  / 12700.0
}

Set-ShapePos $sProblem.Shapes.Item(3) 6606597 3879851    # Picture 4
Set-ShapePos $sProblem.Shapes.Item(4) 6606597 1825625    # Picture 5
Set-ShapePos $sProblem.Shapes.Item(5) 9421798 1825625    # Picture 6
Set-ShapePos $sProblem.Shapes.Item(6) 7315258 2335371    # Rectangle 7
Set-ShapePos $sProblem.Shapes.Item(7) 7315258 2852737    # Rectangle 8
Set-ShapePos $sProblem.Shapes.Item(8) 10157518 2335371   # Rectangle 9
Set-ShapePos $sProblem.Shapes.Item(9) 10157518 2533491   # Rectangle 10
